$d = $word.ActiveDocument

function Get-ParagraphStartingWith($prefix, $after) {
    $count = $d.Paragraphs.Count
    for ($i = $after + 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

# --- Change 1 -----------------------------------------------------------
# "Petra (neutral skeptical):" -> "Petra (neutral frowning):"
# (the speaker tag right before Petra flags down a waiter for a warm drink)
$i1 = Get-ParagraphStartingWith "Petra (neutral skeptical):" 0
$p1 = $d.Paragraphs.Item($i1)
$p1.Range.Find.Execute(
    "Petra (neutral skeptical):", $true, $false, $false, $false, $false,
    $true, 1, $false, "Petra (neutral frowning):", 2) | Out-Null

# --- Change 2 -------------------------------------------------------------
# "Petra (neutral raised_eyebrow): So? ..." -> "Petra (neutral skeptical): So? ..."
$i2 = Get-ParagraphStartingWith "Petra (neutral raised_eyebrow):" $i1
$p2 = $d.Paragraphs.Item($i2)
$p2.Range.Find.Execute(
    "Petra (neutral raised_eyebrow):", $true, $false, $false, $false, $false,
    $true, 1, $false, "Petra (neutral skeptical):", 2) | Out-Null

# The very next paragraph just repeats the bare "Petra (neutral skeptical):"
# speaker tag on its own line -- it gets folded away entirely now that the
# tag lives on the previous line with the dialogue.
$i3 = $i2 + 1
$p3 = $d.Paragraphs.Item($i3)
if ($p3.Range.Text.StartsWith("Petra (neutral skeptical):")) {
    $p3.Range.Delete()
}

# --- Change 3 ---------------------------------------------------------------
# Merge the three runs "I " + "explain" + " to her everything ..." into a
# single run of text (same combined wording, just no longer split across
# separate runs).
$i4 = Get-ParagraphStartingWith "I explain to her everything that happened earlier today" 0
$p4 = $d.Paragraphs.Item($i4)
$r = $p4.Range
$r.MoveEnd(1, -1) | Out-Null   # exclude the trailing paragraph mark
$full = $r.Text
$r.Find.Execute($full, $true, $false, $false, $false, $false, `
    $true, 1, $false, $full, 2) | Out-Null
